$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.557.35"
$ws.Range("E2").Value = "  +5.39%  "
$ws.Range("D3").Value = "1.723.95"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.09"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5387"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06612"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07721"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D13").Value = "1.716.03"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "1.963.13"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5877"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "0.0₅8315"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.00"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").Value = "27.585.44"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.79"
$ws.Range("E19").Value = "  +15.51%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.738"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.102"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.24"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.697"
$ws.Range("E26").Value = "  +12.09%  "
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.66"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05547"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.303"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.548"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.457"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.662"
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9638"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.820"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5961"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01648"
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.928"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.057.59"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8541"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.43"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "1.868.27"
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("E46").Value = "  +12.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.09"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.202"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05280"
$ws.Range("E51").Value = "  +1.80%  "
